$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "62.462.96"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +9.42%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.373.31"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.88%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "415.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.82%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "117.16"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +8.35%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.365.61"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.75%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.577"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.06%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.632"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.15%  "

$ws.Range("E11").Value = "  +19.41%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "40.15"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.28%  "

$ws.Range("E13").Value = "  -0.53%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.895.91"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.59%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "8.35"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.23%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "19.38"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.364.82"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.49%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "62.203.40"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +9.20%  "

$ws.Range("E19").Value = "  -1.98%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.88"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("E21").Value = "  +7.98%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.36"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "12.58"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.82%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "297.12"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "74.92"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  -0.74%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "29.58"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.33%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.05"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +11.75%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.175"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.09%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.24"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.79%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "43.12"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +8.18%  "

$ws.Range("E33").Value = "  +4.56%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "11.46"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.25%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.55"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +19.99%  "

$ws.Range("E36").Value = "  -0.05%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0490"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "52.37"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.58%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +5.91%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.45"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.65%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "133.48"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.54%  "

$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("E44").Value = "  +0.53%  "

$ws.Range("E45").Value = "  +2.80%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.56%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "16.48"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.44%  "

$ws.Range("E48").Value = "  -3.81%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.172.10"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "21.19"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.56%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.692.97"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.42%  "
